$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
# Force text format on Price column cells that hold numeric-looking strings
# so Excel does not coerce them (e.g. "11.50" -> 11.5, "1.001" -> 1.001 as number).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.444.87"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "1.908.94"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "327.42"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").Value = "0.4078"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "0.08015"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "1.007"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "22.33"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").Value = "1.904.83"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "5.936"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "89.09"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "0.06592"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "17.74"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "29.459.65"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "5.536"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "11.50"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").Value = "2.210"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "2.124.65"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "153.33"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").Value = "19.77"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D31").Value = "116.79"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  +9.61%  "
$ws.Range("D33").Value = "0.09474"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "3.579"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").Value = "5.382"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "0.02253"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.06074"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "8.362"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "0.5864"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("D42").Value = "0.1834"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "10.10"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "1.301"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").Value = "0.07748"
$ws.Range("E45").Value = "  +10.25%  "
$ws.Range("D46").Value = "2.375"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").Value = "0.5544"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "12.12"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "1.924"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").Value = "113.34"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "0.2930"
$ws.Range("E51").Value = "  +5.20%  "

# --- Row 29 / 30: ranking order swapped between InternetComputer(DFINITY) and LidoDAOToken ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.127"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.698"
$ws.Range("E30").Value = "  +5.72%  "
